# daily auto push: 2026-01-26 22:37 UTC
#
# A new observation row is inserted right after the existing
# "2026/01/27" entry (row 731), pushing every subsequent row down by
# one. The new row carries: date 2026/01/27, weekday 火, time 5,
# ranking 195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 732:end down by one to make room for the new record.
$ws.Rows.Item(732).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/01/27"), not
# real Excel dates. Force the cell to text first so the slash-delimited
# string isn't auto-coerced into a date serial, then drop the
# formatting override so the cell ends up unstyled like its neighbours.
$ws.Range("A732").NumberFormat = "@"
$ws.Range("A732").Value = "2026/01/27"
$ws.Range("A732").ClearFormats()

$ws.Range("B732").Value = "火"
$ws.Range("C732").Value = 5
$ws.Range("D732").Value = 195
